$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Array / mediani di due array ordinati
$ws.Range("A4").Value = "Array"
$ws.Range("B4").Value = "mediani di due array ordinati"
$ws.Range("C4").Value = "creo la somma dei due array, la sorto e analizzo la sua dimensione: se pari, prendo i valori ris[len(ris)//2] e ris[len(ris//2)-1 e ne faccio la media; se dispari, prendo solo il primo"
$ws.Range("D4").Value = "O(log (n+m)"
$ws.Range("E4").Value = "difficile"

# Row 5: Linked List / merge di due linked list ordinate
$ws.Range("A5").Value = "Linked List"
$ws.Range("B5").Value = "merg di due linked list ordinate"
$ws.Range("C5").Value = "itero in parallelo in entrambe e metto il valore minimo"
$ws.Range("D5").Value = "O(min(m,n))\"
$ws.Range("E5").Value = "facile"

$ws.Range("E5").Select()
